# Refresh the cryptos list (price + 1h volume change columns) with the
# latest scraped figures. Price-like cells that look like plain decimal
# numbers ("42.74", "0.999", ...) need to stay TEXT (the sheet stores
# them as strings, not numbers), so for those we force NumberFormat to
# "@" before the assignment and then restore the cell's prior Style so
# no visible formatting changes.
function Set-TextValue($range, $value) {
    $oldStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $oldStyle
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue $ws.Range("D2") "42.887.97"
$ws.Range("E2").Value = "  +0.18%  "

Set-TextValue $ws.Range("D3") "2.535.85"
$ws.Range("E3").Value = "  -0.85%  "

Set-TextValue $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  -0.06%  "

Set-TextValue $ws.Range("D5") "312.45"
$ws.Range("E5").Value = "  +0.49%  "

Set-TextValue $ws.Range("D6") "100.42"
$ws.Range("E6").Value = "  +2.12%  "

Set-TextValue $ws.Range("D7") "0.566"
$ws.Range("E7").Value = "  -0.77%  "

Set-TextValue $ws.Range("D9") "0.524"
$ws.Range("E9").Value = "  -1.26%  "

Set-TextValue $ws.Range("D10") "35.79"
$ws.Range("E10").Value = "  +0.55%  "

Set-TextValue $ws.Range("D11") "0.0808"
$ws.Range("E11").Value = "  -0.06%  "

Set-TextValue $ws.Range("D12") "7.34"
$ws.Range("E12").Value = "  -1.09%  "

$ws.Range("E13").Value = "  +1.52%  "

Set-TextValue $ws.Range("D14") "2.922.26"
$ws.Range("E14").Value = "  -1.01%  "

Set-TextValue $ws.Range("D15") "15.46"
$ws.Range("E15").Value = "  -3.30%  "

Set-TextValue $ws.Range("D16") "2.522.60"
$ws.Range("E16").Value = "  -3.16%  "

Set-TextValue $ws.Range("D17") "0.820"
$ws.Range("E17").Value = "  -2.16%  "

Set-TextValue $ws.Range("D18") "42.856.67"
$ws.Range("E18").Value = "  +0.13%  "

Set-TextValue $ws.Range("D19") "6.69"
$ws.Range("E19").Value = "  -0.72%  "

Set-TextValue $ws.Range("D20") "12.45"
$ws.Range("E20").Value = "  +1.01%  "

$ws.Range("E21").Value = "  -0.44%  "

Set-TextValue $ws.Range("D22") "69.84"
$ws.Range("E22").Value = "  +0.56%  "

Set-TextValue $ws.Range("D23") "244.26"
$ws.Range("E23").Value = "  -1.33%  "

$ws.Range("E24").Value = "  -0.90%  "

$ws.Range("E25").Value = "  +0.08%  "

$ws.Range("E26").Value = "  +0.05%  "

Set-TextValue $ws.Range("D27") "25.69"
$ws.Range("E27").Value = "  -4.22%  "

$ws.Range("E28").Value = "  -2.32%  "

Set-TextValue $ws.Range("D29") "10.26"
$ws.Range("E29").Value = "  +0.91%  "

Set-TextValue $ws.Range("D30") "38.97"
$ws.Range("E30").Value = "  -2.31%  "

Set-TextValue $ws.Range("D31") "160.60"
$ws.Range("E31").Value = "  +1.57%  "

Set-TextValue $ws.Range("D32") "5.85"
$ws.Range("E32").Value = "  +1.72%  "

$ws.Range("E33").Value = "  +7.77%  "

Set-TextValue $ws.Range("D34") "0.0793"
$ws.Range("E34").Value = "  -0.38%  "

Set-TextValue $ws.Range("D35") "2.67"
$ws.Range("E35").Value = "  +0.99%  "

Set-TextValue $ws.Range("D36") "18.49"
$ws.Range("E36").Value = "  -1.12%  "

Set-TextValue $ws.Range("D37") "3.16"
$ws.Range("E37").Value = "  -3.94%  "

Set-TextValue $ws.Range("D38") "1.98"
$ws.Range("E38").Value = "  -5.12%  "

$ws.Range("E39").Value = "  -0.01%  "

Set-TextValue $ws.Range("D40") "0.118"
$ws.Range("E40").Value = "  +0.26%  "

Set-TextValue $ws.Range("D41") "4.20"
$ws.Range("E41").Value = "  +3.61%  "

Set-TextValue $ws.Range("D42") "22.04"
$ws.Range("E42").Value = "  -2.85%  "

Set-TextValue $ws.Range("D43") "3.36"
$ws.Range("E43").Value = "  +4.93%  "

Set-TextValue $ws.Range("D45") "0.0300"
$ws.Range("E45").Value = "  -0.45%  "

Set-TextValue $ws.Range("D46") "2.004.60"
$ws.Range("E46").Value = "  +0.87%  "

Set-TextValue $ws.Range("D47") "9.28"
$ws.Range("E47").Value = "  +3.55%  "

Set-TextValue $ws.Range("D48") "2.776.56"
$ws.Range("E48").Value = "  -1.00%  "

$ws.Range("E49").Value = "  -0.71%  "

Set-TextValue $ws.Range("D50") "79.89"
$ws.Range("E50").Value = "  -1.83%  "

Set-TextValue $ws.Range("D51") "72.57"
$ws.Range("E51").Value = "  -1.06%  "
